$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 5.828753
$ws.Range("H2").Value = 17.486259
$ws.Range("I2").Value = 0.1911291943607339
$ws.Range("J2").Value = 0.1911291943607339
$ws.Range("M2").Value = 0.01393633333333333
$ws.Range("N2").Value = 0.041809
$ws.Range("Q2").Value = 0.08123144472566667
$ws.Range("R2").Value = 0.731083002531
$ws.Range("S2").Value = 0.1911291943607339
$ws.Range("T2").Value = 0.1911291943607339

# Row 3
$ws.Range("D3").Value = "MuSCs"
$ws.Range("I3").Value = 0.7732994524709527
$ws.Range("J3").Value = 0.7732994524709526
$ws.Range("M3").Value = 0.01393633333333333
$ws.Range("N3").Value = 0.041809
$ws.Range("Q3").Value = 0.3286584864226667
$ws.Range("R3").Value = 2.957926377804
$ws.Range("S3").Value = 0.7732994524709527
$ws.Range("T3").Value = 0.7732994524709526

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 1.084798333333333
$ws.Range("H4").Value = 3.254395
$ws.Range("I4").Value = 0.03557135316831352
$ws.Range("J4").Value = 0.03557135316831351
$ws.Range("M4").Value = 0.01393633333333333
$ws.Range("N4").Value = 0.041809
$ws.Range("Q4").Value = 0.01511811117277778
$ws.Range("R4").Value = 0.136063000555
$ws.Range("S4").Value = 0.03557135316831352
$ws.Range("T4").Value = 0.03557135316831351
